# Add "Area" / "Atotal" columns (G/H) to the discharge sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Row 2: first area segment uses 0 as the baseline depth, plus the running total.
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Row 3: second area segment (not part of the shared-formula fill series below).
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 share the same relative formula pattern (D[r]-D[r-1])*B[r]/100.
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

$ws.Range("D13").Select() | Out-Null

$wb.Save()
